# Update Shimron Hetmyer's Delhi Capitals innings-by-innings stats
# (runs/balls/fours/sixes for rows 2-11) to the corrected figures.
# Values are written with a leading apostrophe so Excel keeps storing
# them as text (matching the original t="str" cell type) instead of
# silently re-typing the numeric-looking strings as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "7"
    "E2"  = "0"
    "F2"  = "1"
    "C3"  = "7"
    "D3"  = "5"
    "E3"  = "0"
    "F3"  = "1"
    "C4"  = "5"
    "D4"  = "5"
    "E4"  = "1"
    "F4"  = "0"
    "C5"  = "42"
    "D5"  = "22"
    "E5"  = "4"
    "C6"  = "11"
    "D6"  = "13"
    "E6"  = "1"
    "F6"  = "0"
    "C7"  = "10"
    "D7"  = "6"
    "E7"  = "0"
    "F7"  = "1"
    "D9"  = "13"
    "E9"  = "1"
    "F9"  = "0"
    "C10" = "16"
    "E10" = "3"
    "C11" = "10"
    "D11" = "5"
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = "'" + $updates[$cellRef]
}
